$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 8 (everything from old row 8 down shifts to row 9+)
$ws.Rows(8).Insert()

# Fill in the new row's data
$ws.Range("A8").Value = "Energia.png"
$ws.Range("B8").Value = "https://opengameart.org/content/energy-icon"
$ws.Range("C8").Value = "CC-BY 4.0"

# Match the "Hyperlink" cell style used by the other Source cells (copy format
# from an existing hyperlink cell so we reuse the same style record instead of
# minting a new one)
$ws.Range("B2").Copy()
$ws.Range("B8").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Rebuild the hyperlinks collection in the final row order so the relationship
# ids line up the same way Excel would after the row insert
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B6"), "https://opengameart.org/content/spikes-0")
$ws.Hyperlinks.Add($ws.Range("B4"), "https://www.fontspace.com/a-area-kilometer-50-font-f53888")
$ws.Hyperlinks.Add($ws.Range("B2"), "https://freesound.org/people/Whiprealgood/sounds/87535/")
$ws.Hyperlinks.Add($ws.Range("B3"), "https://freesound.org/people/suntemple/sounds/253172/")
$ws.Hyperlinks.Add($ws.Range("B5"), "https://opengameart.org/content/simple-explosion-bleeds-game-art")
$ws.Hyperlinks.Add($ws.Range("B7"), "https://opengameart.org/content/various-inventory-24-pixel-icon-set")
$ws.Hyperlinks.Add($ws.Range("B10"), "https://elthen.itch.io/2d-pixel-art-vegetable-monsters-sprite-pack")
$ws.Hyperlinks.Add($ws.Range("B11"), "https://free-game-assets.itch.io/night-city-street-2d-background-tiles")
$ws.Hyperlinks.Add($ws.Range("B8"), "https://opengameart.org/content/energy-icon")

# Reflect the saved selection from the source edit
$ws.Range("C8").Select()
